# Updated legacy GSC export data:
#   - The oldest day in the rolling window (2025-10-23) is dropped.
#   - All remaining days shift up by one row.
#   - The newest day (2026-01-21) is appended at the bottom with 0/0 counts.
#
# This is the "Chart" worksheet that holds the Date / Non-HTTPS URLs / HTTPS URLs
# table (row 1 = headers, rows 2..N = one row per day).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Drop the oldest day's row (row 2, 2025-10-23). Excel shifts every row below it
# up by one, which automatically moves all the other dates/counts into their new
# (one row higher) position - matching the diff's C2/C3/C4 value shifts.
$ws.Rows.Item(2).Delete() | Out-Null

# Figure out the new last row (used range now ends one row earlier) and append
# the newest day after it.
$newRowIndex = $ws.UsedRange.Rows.Count + 1
$newRowRange = "A" + $newRowIndex

# Write the new date as literal text (not an auto-converted date serial number)
# by entering it as a formula returning a text literal, then collapsing that
# formula down to a static value via copy / paste-values. This preserves the
# cell's existing (default) style/number format.
$ws.Range($newRowRange).Formula = '="2026-01-21"'
$ws.Range($newRowRange).Copy() | Out-Null
$ws.Range($newRowRange).PasteSpecial(-4163) | Out-Null   # xlPasteValues
$excel.CutCopyMode = 0

# New day has no recorded URLs yet, consistent with the existing trailing rows.
$ws.Range("B" + $newRowIndex).Value = 0.0
$ws.Range("C" + $newRowIndex).Value = 0.0
